$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 82.38461
$ws.Range("I11").Value = 82.38461
$ws.Range("K11").Value = 82.38461
$ws.Range("M11").Value = 57.61539
$ws.Range("H38").Value = 2830.4
$ws.Range("I38").Value = 49.142857
$ws.Range("J38").Value = 9320
$ws.Range("K38").Value = 147.428571
$ws.Range("L38").Value = 27960
$ws.Range("M38").Value = 224.571429
$ws.Range("N38").Value = -28704
$ws.Range("H41").Value = 653.3
$ws.Range("I41").Value = 653.3
$ws.Range("K41").Value = 653.3
$ws.Range("M41").Value = -213.3
$ws.Range("H88").Value = 1730.7142
$ws.Range("I88").Value = 1598.3334
$ws.Range("J88").Value = 1830
$ws.Range("K88").Value = 1598.3334
$ws.Range("L88").Value = 1830
$ws.Range("M88").Value = -1192.3334
$ws.Range("N88").Value = -2642
$ws.Range("H91").Value = 1730.7142
$ws.Range("I91").Value = 1598.3334
$ws.Range("J91").Value = 1830
$ws.Range("K91").Value = 1598.3334
$ws.Range("L91").Value = 1830
$ws.Range("M91").Value = -194.3334
$ws.Range("N91").Value = -4638
$ws.Range("H92").Value = 203.11111
$ws.Range("I92").Value = 212.5
$ws.Range("J92").Value = 195.6
$ws.Range("K92").Value = 212.5
$ws.Range("L92").Value = 195.6
$ws.Range("M92").Value = 1035.5
$ws.Range("N92").Value = -2691.6
$ws.Range("H125").Value = 1065.6666
$ws.Range("J125").Value = 550
$ws.Range("L125").Value = 4950
$ws.Range("N125").Value = -9870
$ws.Range("H132").Value = 2917.05
$ws.Range("I132").Value = 1685.6111
$ws.Range("K132").Value = 5056.8333
$ws.Range("M132").Value = -2526.8333
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1741.6666
$ws.Range("I88").Value = 1742.5
$ws.Range("K88").Value = 1742.5
$ws.Range("M88").Value = -1336.5
$ws.Range("H91").Value = 1741.6666
$ws.Range("I91").Value = 1742.5
$ws.Range("K91").Value = 1742.5
$ws.Range("M91").Value = -338.5
$ws.Range("H102").Value = 10421565
$ws.Range("I102").Value = 20833576
$ws.Range("K102").Value = 20833576
$ws.Range("M102").Value = -20831954
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H122").Value = 800
$ws.Range("I122").Value = 800
$ws.Range("K122").Value = 2400
$ws.Range("M122").Value = 50
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2898.4285
$ws.Range("I86").Value = 1324.5333
$ws.Range("J86").Value = 6833.1665
$ws.Range("K86").Value = 1324.5333
$ws.Range("L86").Value = 6833.1665
$ws.Range("M86").Value = -201.5333000000001
$ws.Range("N86").Value = -9079.166499999999
$ws.Range("H89").Value = 2898.4285
$ws.Range("I89").Value = 1324.5333
$ws.Range("J89").Value = 6833.1665
$ws.Range("K89").Value = 6622.6665
$ws.Range("L89").Value = 34165.8325
$ws.Range("M89").Value = -1006.6665
$ws.Range("N89").Value = -45397.8325
$ws.Range("H105").Value = 3953818.5
$ws.Range("I105").Value = 5051525
$ws.Range("K105").Value = 5051525
$ws.Range("M105").Value = -5049778
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4384.778
$ws.Range("I31").Value = 1081.7059
$ws.Range("J31").Value = 10000
$ws.Range("K31").Value = 1081.7059
$ws.Range("L31").Value = 10000
$ws.Range("M31").Value = -786.7058999999999
$ws.Range("N31").Value = -10590
$ws.Range("H33").Value = 1500
$ws.Range("I33").Value = 1500
$ws.Range("K33").Value = 1500
$ws.Range("M33").Value = -1121
$ws.Range("H34").Value = 4384.778
$ws.Range("I34").Value = 1081.7059
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 1081.7059
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = -879.7058999999999
$ws.Range("N34").Value = -10404
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31996
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99984
$ws.Range("H88").Value = 7823.25
$ws.Range("J88").Value = 7823.25
$ws.Range("L88").Value = 7823.25
$ws.Range("N88").Value = -8635.25
$ws.Range("H91").Value = 7823.25
$ws.Range("J91").Value = 7823.25
$ws.Range("L91").Value = 7823.25
$ws.Range("N91").Value = -10631.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1402.1177
$ws.Range("J5").Value = 1409.1428
$ws.Range("L5").Value = 4227.428400000001
$ws.Range("N5").Value = -4451.428400000001
$ws.Range("H22").Value = 216.66667
$ws.Range("H27").Value = 216.66667
$ws.Range("H44").Value = 3561
$ws.Range("J44").Value = 10300
$ws.Range("L44").Value = 30900
$ws.Range("N44").Value = -31696
$ws.Range("H135").Value = 1402.1177
$ws.Range("J135").Value = 1409.1428
$ws.Range("L135").Value = 12682.2852
$ws.Range("N135").Value = -17752.2852
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2125
$ws.Range("I80").Value = 2125
$ws.Range("K80").Value = 2125
$ws.Range("M80").Value = -1127
$ws.Range("H83").Value = 2125
$ws.Range("I83").Value = 2125
$ws.Range("K83").Value = 10625
$ws.Range("M83").Value = -5633
$ws.Range("H122").Value = 2089.5
$ws.Range("I122").Value = 1493.2222
$ws.Range("K122").Value = 4479.6666
$ws.Range("M122").Value = -2029.6666
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1214.6923
$ws.Range("I22").Value = 978.4
$ws.Range("J22").Value = 1362.375
$ws.Range("K22").Value = 978.4
$ws.Range("L22").Value = 1362.375
$ws.Range("M22").Value = -683.4
$ws.Range("N22").Value = -1952.375
$ws.Range("H27").Value = 1214.6923
$ws.Range("I27").Value = 978.4
$ws.Range("J27").Value = 1362.375
$ws.Range("K27").Value = 978.4
$ws.Range("L27").Value = 1362.375
$ws.Range("M27").Value = -871.4
$ws.Range("N27").Value = -1576.375
$ws.Range("H40").Value = 1884.125
$ws.Range("I40").Value = 1884.125
$ws.Range("K40").Value = 1884.125
$ws.Range("M40").Value = -1748.125
$ws.Range("H42").Value = 30012
$ws.Range("I42").Value = 20025
$ws.Range("K42").Value = 20025
$ws.Range("M42").Value = -19462
$ws.Range("H49").Value = 30012
$ws.Range("I49").Value = 20025
$ws.Range("K49").Value = 20025
$ws.Range("M49").Value = -19878
$ws.Range("H82").Value = 2744.4119
$ws.Range("I82").Value = 575.875
$ws.Range("K82").Value = 575.875
$ws.Range("M82").Value = -214.875
$ws.Range("H85").Value = 2744.4119
$ws.Range("I85").Value = 575.875
$ws.Range("K85").Value = 575.875
$ws.Range("M85").Value = 672.125
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 316.33334
$ws.Range("J6").Value = 316.33334
$ws.Range("L6").Value = 316.33334
$ws.Range("N6").Value = -546.33334
$ws.Range("H107").Value = 37037860
$ws.Range("I107").Value = 55556136
$ws.Range("J107").Value = 1312
$ws.Range("K107").Value = 166668408
$ws.Range("L107").Value = 3936
$ws.Range("M107").Value = -166666488
$ws.Range("N107").Value = -7776
$ws.Range("H122").Value = 4951
$ws.Range("I122").Value = 2902
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 8706
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -6256
$ws.Range("N122").Value = -25900
$ws.Range("H126").Value = 3768
$ws.Range("I126").Value = 1614.3636
$ws.Range("K126").Value = 4843.0908
$ws.Range("M126").Value = -2373.0908
$ws.Range("H136").Value = 2944
$ws.Range("I136").Value = 1821.3636
$ws.Range("K136").Value = 5464.0908
$ws.Range("M136").Value = -2914.0908
